$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data for rows 2-57: column A = Elapsed Time (s), column B = Enemy in Sight
$data = @(
    @(1, $false),
    @(2, $false),
    @(3, $false),
    @(4, $false),
    @(5, $false),
    @(5, $false),
    @(6, $false),
    @(7, $false),
    @(8, $true),
    @(9, $false),
    @(10, $false),
    @(11, $false),
    @(12, $false),
    @(12, $false),
    @(13, $true),
    @(14, $false),
    @(15, $false),
    @(16, $false),
    @(17, $false),
    @(18, $false),
    @(19, $false),
    @(20, $false),
    @(20, $false),
    @(21, $false),
    @(22, $false),
    @(23, $false),
    @(24, $false),
    @(25, $false),
    @(26, $false),
    @(26, $false),
    @(27, $false),
    @(28, $false),
    @(29, $false),
    @(30, $false),
    @(31, $false),
    @(32, $false),
    @(33, $false),
    @(33, $false),
    @(34, $false),
    @(35, $false),
    @(36, $false),
    @(37, $false),
    @(38, $false),
    @(39, $false),
    @(40, $true),
    @(41, $true),
    @(41, $false),
    @(42, $false),
    @(43, $false),
    @(44, $true),
    @(45, $false),
    @(46, $false),
    @(47, $false),
    @(48, $false),
    @(48, $false),
    @(49, $false)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
}
